$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 647, shifting existing rows 647:736 down to 648:737.
$ws.Rows.Item(647).Insert()

# Populate the newly inserted row 647 with the new record's data.
$ws.Cells.Item(647, 1).Value = 3
$ws.Cells.Item(647, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(647, 3).Value = "Coquimbo"
$ws.Cells.Item(647, 4).Value = Get-Date -Year 2023 -Month 7 -Day 24 -Hour 0 -Minute 0 -Second 0
$ws.Cells.Item(647, 5).Value = 5
$ws.Cells.Item(647, 6).Value = 100112021
$ws.Cells.Item(647, 7).Value = "Ají"
$ws.Cells.Item(647, 8).Value = "Inferno"
$ws.Cells.Item(647, 9).Value = "Primera"
$ws.Cells.Item(647, 10).Value = 75
$ws.Cells.Item(647, 11).Value = 12000
$ws.Cells.Item(647, 12).Value = 12500
$ws.Cells.Item(647, 13).Value = 12267
$ws.Cells.Item(647, 14).Value = "$/caja 10 kilos"
$ws.Cells.Item(647, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(647, 16).Value = 1227
$ws.Cells.Item(647, 17).Value = 10
$ws.Cells.Item(647, 18).Value = "Hortaliza"
